$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Previously added")
$ws2 = $wb.Worksheets.Item("New")

function Set-TextValue($ws, $row, $col, $val) {
    # Force the cell to be stored as a shared string even when the text looks
    # numeric (e.g. cadastre numbers like "40460060490") by switching the
    # cell to a text number format before assigning the value.
    $ws.Cells.Item($row, $col).NumberFormat = "@"
    $ws.Cells.Item($row, $col).Value = $val
}

function Restore-Format($wsSrc, $srcCell, $wsDest, $destCell) {
    # Copy ONLY formatting (styles) from a known-good, already-saved data
    # cell so the destination ends up with the same cellXfs index (s=3 for
    # links, s=4 for plain text, s=2 for the date column) instead of
    # whatever ad-hoc style got stamped on it by NumberFormat/Hyperlinks.Add.
    $wsSrc.Range($srcCell).Copy()
    $wsDest.Range($destCell).PasteSpecial(-4122)
}

# =====================================================================
# Step 1: archive the three current "New" sheet listings (rows 2-4) onto
# the end of "Previously added" (rows 512-514), exactly as they stand now.
# =====================================================================

$archiveRows = @(
    @("https://www.ss.com/msg/lv/real-estate/wood/bauska-and-reg/iecavas-nov/lclfm.html", "7 000 €", "Bauska un raj.", "4 ha.", "40460060490", 46071.69236111111),
    @("https://www.ss.com/msg/lv/real-estate/wood/kraslava-and-reg/kepovas-pag/mxnne.html", "65 000 €", "Krāslava un raj.", "9.50 ha.", "60800040007", 46072.018055555556),
    @("https://www.ss.com/msg/lv/real-estate/wood/preili-and-reg/rozupes-pag/ilhhk.html", "10 €", "Preiļi un raj.", "3 ha.", "76660010146", 46072.54166666667)
)

$destStart = 512
for ($i = 0; $i -le 2; $i++) {
    $destRow = $destStart + $i
    $vals = $archiveRows[$i]

    Set-TextValue $ws1 $destRow 1 $vals[0]
    Set-TextValue $ws1 $destRow 2 $vals[1]
    Set-TextValue $ws1 $destRow 3 $vals[2]
    Set-TextValue $ws1 $destRow 4 $vals[3]
    Set-TextValue $ws1 $destRow 5 $vals[4]
    $ws1.Cells.Item($destRow, 6).Value = $vals[5]

    $ws1.Hyperlinks.Add($ws1.Cells.Item($destRow, 1), $vals[0])

    # re-apply the standard data-row formatting (styles s=3/4/2) last, since
    # both the text-number-format trick and Hyperlinks.Add stamp their own
    # styles on the cell. Row 511 (the last untouched data row) is the
    # known-good style template.
    Restore-Format $ws1 "A511" $ws1 ("A" + $destRow)
    Restore-Format $ws1 "B511" $ws1 ("B" + $destRow)
    Restore-Format $ws1 "C511" $ws1 ("C" + $destRow)
    Restore-Format $ws1 "D511" $ws1 ("D" + $destRow)
    Restore-Format $ws1 "E511" $ws1 ("E" + $destRow)
    Restore-Format $ws1 "F511" $ws1 ("F" + $destRow)
}
$excel.CutCopyMode = 0

# =====================================================================
# Step 2: replace "New" sheet rows 2-4 with the freshly scraped listings.
# =====================================================================

$newRows = @(
    @("https://www.ss.com/msg/lv/real-estate/wood/balvi-and-reg/viksnas-pag/eeedd.html", "20 000 €", "Balvi un raj.", "2.27 ha.", "38940060037", 46072.90694444445),
    @("https://www.ss.com/msg/lv/real-estate/wood/gulbene-and-reg/stradu-pag/nlppl.html", "40 000 €", "Gulbene un raj.", "12 ha.", "", 46073.47083333333),
    @("https://www.ss.com/msg/lv/real-estate/wood/preili-and-reg/preilu-pag/cckmx.html", "5 000 €", "Preiļi un raj.", "3 ha.", "76660010146", 46072.83888888889)
)

for ($i = 0; $i -le 2; $i++) {
    $row = 2 + $i
    $vals = $newRows[$i]

    Set-TextValue $ws2 $row 1 $vals[0]
    Set-TextValue $ws2 $row 2 $vals[1]
    Set-TextValue $ws2 $row 3 $vals[2]
    Set-TextValue $ws2 $row 4 $vals[3]
    if ($vals[4] -ne "") {
        Set-TextValue $ws2 $row 5 $vals[4]
    } else {
        $ws2.Cells.Item($row, 5).Value = ""
    }
    $ws2.Cells.Item($row, 6).Value = $vals[5]

    # same style clean-up as above, sourced from the same reference row.
    Restore-Format $ws1 "A511" $ws2 ("A" + $row)
    Restore-Format $ws1 "B511" $ws2 ("B" + $row)
    Restore-Format $ws1 "C511" $ws2 ("C" + $row)
    Restore-Format $ws1 "D511" $ws2 ("D" + $row)
    Restore-Format $ws1 "E511" $ws2 ("E" + $row)
    Restore-Format $ws1 "F511" $ws2 ("F" + $row)
}
$excel.CutCopyMode = 0

# update the existing hyperlinks on "New" rows 2-4 to point at the new URLs
$i = 0
foreach ($h in $ws2.Hyperlinks) {
    $i = $i + 1
    if ($i -ge 1 -and $i -le 3) {
        $h.Address = $newRows[$i - 1][0]
    }
}
